$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.020187983908264
$ws.Cells.Item(2, 4).Value = 1.028019513462228
$ws.Cells.Item(2, 5).Value = 1.021231247711937
$ws.Cells.Item(2, 6).Value = 1.034508892943463
$ws.Cells.Item(2, 9).Value = 1.024714453624212
$ws.Cells.Item(2, 10).Value = 1.025386507826116
$ws.Cells.Item(2, 11).Value = 1.030837548503848
$ws.Cells.Item(2, 12).Value = 1.024069188873236
$ws.Cells.Item(2, 13).Value = 1.037308163757292
$ws.Cells.Item(2, 14).Value = 1.012546891296607

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.021407650989037
$ws.Cells.Item(3, 4).Value = 1.029206442331322
$ws.Cells.Item(3, 5).Value = 1.022270910416079
$ws.Cells.Item(3, 6).Value = 1.036038823618318
$ws.Cells.Item(3, 9).Value = 1.024898766294557
$ws.Cells.Item(3, 10).Value = 1.026241844292395
$ws.Cells.Item(3, 11).Value = 1.031831301572135
$ws.Cells.Item(3, 12).Value = 1.024914621453541
$ws.Cells.Item(3, 13).Value = 1.038645378321752
$ws.Cells.Item(3, 14).Value = 1.012828302403916

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.022191179941
$ws.Cells.Item(4, 4).Value = 1.029966393190839
$ws.Cells.Item(4, 5).Value = 1.022939141686442
$ws.Cells.Item(4, 6).Value = 1.037011838863312
$ws.Cells.Item(4, 9).Value = 1.025007816113666
$ws.Cells.Item(4, 10).Value = 1.02678919751559
$ws.Cells.Item(4, 11).Value = 1.032465712636473
$ws.Cells.Item(4, 12).Value = 1.025456590957934
$ws.Cells.Item(4, 13).Value = 1.039493247414664
$ws.Cells.Item(4, 14).Value = 1.013008375584385

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.022519230798228
$ws.Cells.Item(5, 4).Value = 1.030283961230074
$ws.Cells.Item(5, 5).Value = 1.023219000284593
$ws.Cells.Item(5, 6).Value = 1.037416867051381
$ws.Cells.Item(5, 9).Value = 1.025051222160216
$ws.Cells.Item(5, 10).Value = 1.027017855090951
$ws.Cells.Item(5, 11).Value = 1.032730372904478
$ws.Cells.Item(5, 12).Value = 1.025683229274953
$ws.Cells.Item(5, 13).Value = 1.039845555046187
$ws.Cells.Item(5, 14).Value = 1.013083599114939

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.022574233529033
$ws.Cells.Item(6, 4).Value = 1.03033717046023
$ws.Cells.Item(6, 5).Value = 1.023265927594014
$ws.Cells.Item(6, 6).Value = 1.037484637754806
$ws.Cells.Item(6, 9).Value = 1.02505836740066
$ws.Cells.Item(6, 10).Value = 1.027056163031493
$ws.Cells.Item(6, 11).Value = 1.032774690984006
$ws.Cells.Item(6, 12).Value = 1.025721212449738
$ws.Cells.Item(6, 13).Value = 1.039904467269643
$ws.Cells.Item(6, 14).Value = 1.013096201486147

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022195568637305
$ws.Cells.Item(7, 4).Value = 1.029970644053923
$ws.Cells.Item(7, 5).Value = 1.022942885343491
$ws.Cells.Item(7, 6).Value = 1.037017266654247
$ws.Cells.Item(7, 9).Value = 1.025008405680676
$ws.Cells.Item(7, 10).Value = 1.026792258529592
$ws.Cells.Item(7, 11).Value = 1.032469257058876
$ws.Cells.Item(7, 12).Value = 1.02545962403318
$ws.Cells.Item(7, 13).Value = 1.039497971182823
$ws.Cells.Item(7, 14).Value = 1.013009382602635

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.020601360574872
$ws.Cells.Item(8, 4).Value = 1.028422323513221
$ws.Cells.Item(8, 5).Value = 1.021583545505672
$ws.Cells.Item(8, 6).Value = 1.035029469611875
$ws.Cells.Item(8, 9).Value = 1.024778861800762
$ws.Cells.Item(8, 10).Value = 1.025676844848813
$ws.Cells.Item(8, 11).Value = 1.031175185502794
$ws.Cells.Item(8, 12).Value = 1.024355965849289
$ws.Cells.Item(8, 13).Value = 1.03776370202266
$ws.Cells.Item(8, 14).Value = 1.012642415980874

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.017747988211646
$ws.Cells.Item(9, 4).Value = 1.025631340423472
$ws.Cells.Item(9, 5).Value = 1.019153194845662
$ws.Cells.Item(9, 6).Value = 1.031395373929443
$ws.Cells.Item(9, 9).Value = 1.024295812244497
$ws.Cells.Item(9, 10).Value = 1.023663979802029
$ws.Cells.Item(9, 11).Value = 1.028828117282656
$ws.Cells.Item(9, 12).Value = 1.022371726555801
$ws.Cells.Item(9, 13).Value = 1.034573046119411
$ws.Cells.Item(9, 14).Value = 1.011980118202319

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.015815001574533
$ws.Cells.Item(10, 4).Value = 1.02372731865093
$ws.Cells.Item(10, 5).Value = 1.017508579499424
$ws.Cells.Item(10, 6).Value = 1.028882048800113
$ws.Cells.Item(10, 9).Value = 1.023920390371595
$ws.Cells.Item(10, 10).Value = 1.022289328102427
$ws.Cells.Item(10, 11).Value = 1.027217346872778
$ws.Cells.Item(10, 12).Value = 1.021021573354788
$ws.Cells.Item(10, 13).Value = 1.032353273244774
$ws.Cells.Item(10, 14).Value = 1.01152776530346

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.014970455486736
$ws.Cells.Item(11, 4).Value = 1.022892266142729
$ws.Cells.Item(11, 5).Value = 1.016790456793382
$ws.Cells.Item(11, 6).Value = 1.0277717040162
$ws.Cells.Item(11, 9).Value = 1.023745015550867
$ws.Cells.Item(11, 10).Value = 1.021686104073031
$ws.Cells.Item(11, 11).Value = 1.026508652269806
$ws.Cells.Item(11, 12).Value = 1.020430265352416
$ws.Cells.Item(11, 13).Value = 1.031369581687988
$ws.Cells.Item(11, 14).Value = 1.011329252993937

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.014655596203209
$ws.Cells.Item(12, 4).Value = 1.022580470162378
$ws.Cells.Item(12, 5).Value = 1.016522794949769
$ws.Cells.Item(12, 6).Value = 1.027355906890876
$ws.Cells.Item(12, 9).Value = 1.023677933705406
$ws.Cells.Item(12, 10).Value = 1.02146081913837
$ws.Cells.Item(12, 11).Value = 1.026243699142766
$ws.Cells.Item(12, 12).Value = 1.020209605268778
$ws.Cells.Item(12, 13).Value = 1.031000763507274
$ws.Cells.Item(12, 14).Value = 1.01125511335634

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.014723187338724
$ws.Cells.Item(13, 4).Value = 1.022647425247508
$ws.Cells.Item(13, 5).Value = 1.016580251225801
$ws.Cells.Item(13, 6).Value = 1.027445249949818
$ws.Cells.Item(13, 9).Value = 1.023692411018036
$ws.Cells.Item(13, 10).Value = 1.021509199095768
$ws.Cells.Item(13, 11).Value = 1.026300610422804
$ws.Cells.Item(13, 12).Value = 1.020256984151251
$ws.Cells.Item(13, 13).Value = 1.031080032348527
$ws.Cells.Item(13, 14).Value = 1.011271034922099

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.014944452865071
$ws.Cells.Item(14, 4).Value = 1.022866526211028
$ws.Cells.Item(14, 5).Value = 1.016768350642834
$ws.Cells.Item(14, 6).Value = 1.027737403063597
$ws.Cells.Item(14, 9).Value = 1.02373951020534
$ws.Cells.Item(14, 10).Value = 1.021667506972115
$ws.Cells.Item(14, 11).Value = 1.026486786253958
$ws.Cells.Item(14, 12).Value = 1.020412046476358
$ws.Cells.Item(14, 13).Value = 1.031339165352739
$ws.Cells.Item(14, 14).Value = 1.011323132855769

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015080627803149
$ws.Cells.Item(15, 4).Value = 1.02300130597635
$ws.Cells.Item(15, 5).Value = 1.016884122537073
$ws.Cells.Item(15, 6).Value = 1.027916960664788
$ws.Cells.Item(15, 9).Value = 1.023768272065746
$ws.Cells.Item(15, 10).Value = 1.02176488330652
$ws.Cells.Item(15, 11).Value = 1.026601267594866
$ws.Cells.Item(15, 12).Value = 1.020507449518212
$ws.Cells.Item(15, 13).Value = 1.031498369557067
$ws.Cells.Item(15, 14).Value = 1.01135517846304

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.015870890313602
$ws.Cells.Item(16, 4).Value = 1.023782512683991
$ws.Cells.Item(16, 5).Value = 1.017556111095169
$ws.Cells.Item(16, 6).Value = 1.028955269644768
$ws.Cells.Item(16, 9).Value = 1.023931758215562
$ws.Cells.Item(16, 10).Value = 1.022329192044392
$ws.Cells.Item(16, 11).Value = 1.027264141958601
$ws.Cells.Item(16, 12).Value = 1.021060674211793
$ws.Cells.Item(16, 13).Value = 1.03241807915563
$ws.Cells.Item(16, 14).Value = 1.011540883722472

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.016364563906313
$ws.Cells.Item(17, 4).Value = 1.024269686224758
$ws.Cells.Item(17, 5).Value = 1.017976014310292
$ws.Cells.Item(17, 6).Value = 1.029600632285539
$ws.Cells.Item(17, 9).Value = 1.024030868016969
$ws.Cells.Item(17, 10).Value = 1.022681014784317
$ws.Cells.Item(17, 11).Value = 1.027676923184013
$ws.Cells.Item(17, 12).Value = 1.021405895922166
$ws.Cells.Item(17, 13).Value = 1.032988927279602
$ws.Cells.Item(17, 14).Value = 1.011656660227707

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.016651788499792
$ws.Cells.Item(18, 4).Value = 1.024552825684048
$ws.Cells.Item(18, 5).Value = 1.018220359791342
$ws.Cells.Item(18, 6).Value = 1.02997493673063
$ws.Cells.Item(18, 9).Value = 1.024087441713387
$ws.Cells.Item(18, 10).Value = 1.022885457466794
$ws.Cells.Item(18, 11).Value = 1.027916610771627
$ws.Cells.Item(18, 12).Value = 1.021606614384589
$ws.Cells.Item(18, 13).Value = 1.033319724388744
$ws.Cells.Item(18, 14).Value = 1.011723936400369

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.016749602051776
$ws.Cells.Item(19, 4).Value = 1.024649196630351
$ws.Cells.Item(19, 5).Value = 1.018303577972612
$ws.Cells.Item(19, 6).Value = 1.030102206010069
$ws.Cells.Item(19, 9).Value = 1.024106522762355
$ws.Cells.Item(19, 10).Value = 1.022955037198525
$ws.Cells.Item(19, 11).Value = 1.027998155519955
$ws.Cells.Item(19, 12).Value = 1.021674945569252
$ws.Cells.Item(19, 13).Value = 1.033432151192541
$ws.Cells.Item(19, 14).Value = 1.011746832897001

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.016311672746864
$ws.Cells.Item(20, 4).Value = 1.024217522880034
$ws.Cells.Item(20, 5).Value = 1.017931022478776
$ws.Cells.Item(20, 6).Value = 1.029531611109816
$ws.Cells.Item(20, 9).Value = 1.024020362355287
$ws.Cells.Item(20, 10).Value = 1.022643347276369
$ws.Cells.Item(20, 11).Value = 1.027632747636811
$ws.Cells.Item(20, 12).Value = 1.021368923588962
$ws.Cells.Item(20, 13).Value = 1.032927905365306
$ws.Cells.Item(20, 14).Value = 1.011644264854996

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.014879327818716
$ws.Cells.Item(21, 4).Value = 1.022802051398916
$ws.Cells.Item(21, 5).Value = 1.016712985579151
$ws.Cells.Item(21, 6).Value = 1.027651464590511
$ws.Cells.Item(21, 9).Value = 1.023725694340549
$ws.Cells.Item(21, 10).Value = 1.021620923112159
$ws.Cells.Item(21, 11).Value = 1.026432009604241
$ws.Cells.Item(21, 12).Value = 1.020366412846752
$ws.Cells.Item(21, 13).Value = 1.031262952268796
$ws.Cells.Item(21, 14).Value = 1.011307802501592

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.013972047916141
$ws.Cells.Item(22, 4).Value = 1.021902698748916
$ws.Cells.Item(22, 5).Value = 1.015941829851547
$ws.Cells.Item(22, 6).Value = 1.026449841592734
$ws.Cells.Item(22, 9).Value = 1.023529193492989
$ws.Cells.Item(22, 10).Value = 1.020971011228389
$ws.Cells.Item(22, 11).Value = 1.025667136221135
$ws.Cells.Item(22, 12).Value = 1.019730172169728
$ws.Cells.Item(22, 13).Value = 1.030196251276136
$ws.Cells.Item(22, 14).Value = 1.011093918192645

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01445365795338
$ws.Cells.Item(23, 4).Value = 1.022380362383766
$ws.Cells.Item(23, 5).Value = 1.016351145582116
$ws.Cells.Item(23, 6).Value = 1.027088711094957
$ws.Cells.Item(23, 9).Value = 1.023634432106123
$ws.Cells.Item(23, 10).Value = 1.021316219299909
$ws.Cells.Item(23, 11).Value = 1.026073559844422
$ws.Cells.Item(23, 12).Value = 1.020068023020917
$ws.Cells.Item(23, 13).Value = 1.030763630966198
$ws.Cells.Item(23, 14).Value = 1.011207526146224

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016335574230312
$ws.Cells.Item(24, 4).Value = 1.024241096402753
$ws.Cells.Item(24, 5).Value = 1.01795135413294
$ws.Cells.Item(24, 6).Value = 1.029562805370938
$ws.Cells.Item(24, 9).Value = 1.024025113228273
$ws.Cells.Item(24, 10).Value = 1.022660369978987
$ws.Cells.Item(24, 11).Value = 1.027652712003357
$ws.Cells.Item(24, 12).Value = 1.021385631782669
$ws.Cells.Item(24, 13).Value = 1.032955485241466
$ws.Cells.Item(24, 14).Value = 1.011649866575904

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.01849099112684
$ws.Cells.Item(25, 4).Value = 1.026360423321499
$ws.Cells.Item(25, 5).Value = 1.019785731479789
$ws.Cells.Item(25, 6).Value = 1.032350664966911
$ws.Cells.Item(25, 9).Value = 1.024430056342645
$ws.Cells.Item(25, 10).Value = 1.024190055035852
$ws.Cells.Item(25, 11).Value = 1.029442917045443
$ws.Cells.Item(25, 12).Value = 1.022889454452621
$ws.Cells.Item(25, 13).Value = 1.035414074134527
$ws.Cells.Item(25, 14).Value = 1.012153222464592
